# Updates cryptos list price (column D) and 1h volume change (column E)
# values, matching the scheduled GitHub Actions refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.183.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.856.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7136"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9993"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07744"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3073"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.852.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.225"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7166"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.151.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.859"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "243.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007789"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.105.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.992"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.05%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1597"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.909"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.495"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.314"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.403"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.196"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05183"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.172"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7260"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.674"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("E38").Value = "  -0.72%  "
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.157.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9033"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.142"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.996.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5216"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.765"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.307"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.852"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.12%  "
